$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "wong3"

# Update numeric values per the diff
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 14
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 12

$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 3

$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 22
$ws.Range("E5").Value = 56
$ws.Range("F5").Value = 50

$ws.Range("E6").Value = 55
$ws.Range("F6").Value = 55

$ws.Range("E7").Value = 14
$ws.Range("F7").Value = 14

$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 6
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = 12

$ws.Range("E9").Value = 30
$ws.Range("F9").Value = 30

$ws.Range("B10").Value = 16
$ws.Range("C10").Value = 16
$ws.Range("E10").Value = 47
$ws.Range("F10").Value = 47

$ws.Range("B13").Value = 14
$ws.Range("C13").Value = 14
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 13

$ws.Range("B14").Value = 17
$ws.Range("C14").Value = 17
$ws.Range("E14").Value = 19
$ws.Range("F14").Value = 19

$ws.Range("B15").Value = 11
$ws.Range("C15").Value = 11
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 16

$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 13

$ws.Range("E17").Value = 52
$ws.Range("F17").Value = 52

$ws.Range("E18").Value = 12
$ws.Range("F18").Value = 12

$ws.Range("B21").Value = 19
$ws.Range("C21").Value = 19
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 16

$ws.Range("B22").Value = 12
$ws.Range("C22").Value = 12
$ws.Range("E22").Value = 26
$ws.Range("F22").Value = 24
